# Coffey R00 Master Key - data 4/10-11 added, mild fixes
# Update the "notes" (column W) entries on Sheet1 to reflect new session
# observations (ran unplugged, video issues, connector problems, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update / expand existing notes -----------------------------------
$ws.Range("W2").Value  = "2) Connector fell off implant during set-up, didn't run. Reattached later in the day. 4) chewed through cord midsession. 5) Unplugged midsession, won't stop going at cord. Unplugged hereafter. 7) no video."
$ws.Range("W3").Value  = "7) Unplugged midsession, going after cord. No video. 8) Unplugged whole session. 9) Unplugged mid-session. Unplugged hereafter."
$ws.Range("W6").Value  = "3) chewed through cord mid-session.  4) chewed through cord mid-session. 5) unplugged midsession, won't stop going at cord. Unplugged hereafter. 7) no video."
$ws.Range("W7").Value  = "7) no video. 8) connector fell off presession. 9) ran unplugged, connector reattached after session. "
$ws.Range("W15").Value = "9) connector fell off implant mid-session (reattached later in day). 11) whole exterior of the implant came off (pre-session), dental cement and glue intact. Not plugged in again. "
$ws.Range("W17").Value = "7) focused on jumping/grappling at the ceiling hole/grabbing cord whole session. 8) unplugged mid-session, jumping/cord-grabbing. 9) unplugged mid-session, jumping/cord-grabbing again. Unplugged hereafter. 14) replugged in. 15) came unplugged in the last half hour."

# --- New notes for rows that previously had none ----------------------
$ws.Range("W4").Value  = "14) video started late"
$ws.Range("W5").Value  = "14) video started late"
$ws.Range("W8").Value  = "14) video started late"
$ws.Range("W10").Value = "13) ran unplugged by accident. "
$ws.Range("W11").Value = "never plugged in, cap glued on. "
$ws.Range("W13").Value = "15) came unplugged in the last hour"

# --- Note no longer applicable -----------------------------------------
$ws.Range("W14").Clear()

# --- Stray formatted (empty) cell left behind while editing near S21 --
$ws.Range("W9").Copy() | Out-Null
$ws.Range("S21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Widen notes column to fit the newly expanded text -----------------
$ws.Columns.Item(23).ColumnWidth = 204

# --- Restore the active selection to W3 --------------------------------
$ws.Range("W3").Select()
